# Fill in circulation / ILL statistics for July 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(43984,6000,7124),
    @(20642,2194,2499),
    @(72486,4947,6352),
    @(1593,440,55),
    @(42467,7277,6350),
    @(4817,1039,1000),
    @(5651,961,748),
    @(2215,279,86),
    @(388,235,4),
    @(0,0,0),
    @(978,289,222),
    @(2852,1385,890),
    @(5754,1966,963),
    @(3425,1427,597),
    @(1958,623,138),
    @(16993,2503,2982),
    @(2766,793,503),
    @(20706,2138,3522),
    @(296,324,12),
    @(17646,2048,2710),
    @(1133,452,175),
    @(17674,3064,3056),
    @(66345,7699,7378),
    @(5114,1574,663),
    @(0,0,0),
    @(5595,1045,1128),
    @(2545,309,498),
    @(15471,2719,2377),
    @(593,239,216),
    @(2425,1552,282),
    @(14846,2936,2692),
    @(11212,2545,2725),
    @(6366,737,1369),
    @(51734,5382,5295),
    @(8368,2549,1337),
    @(19315,1860,2814),
    @(933,738,137),
    @(1643,257,765),
    @(1827,234,130),
    @(7719,248,192),
    @(240,102,72),
    @(551,48,22),
    @(1743,159,73),
    @(3232,898,466),
    @(11950,3088,2333),
    @(30096,2968,4443),
    @(15085,3088,1098),
    @(10488,1037,1703),
    @(26583,2841,3179),
    @(3721,425,806),
    @(13158,2711,2325),
    @(2160,656,1088),
    @(1820,1316,200),
    @(3406,990,1166),
    @(13009,5067,2672),
    @(12080,791,489),
    @(632708,93334,88960)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 2).Value = $rowVals[0]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
}
